# Insert a new weekly record before the current row 240, pushing the
# existing rows 240-299 down to 241-300 (same as the author's diff, which
# shows every row from 240 to 299 taking on the values previously held by
# the row above it, and a brand-new row 300 appearing with the data that
# used to live in row 299).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 240:299 down to 241:300, formatting included.
$ws.Range("A240").EntireRow.Insert()

# Populate the newly freed row 240 with the new weekly observation.
$ws.Range("A240").Value = 7
$ws.Range("B240").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C240").Value = "Ñuble"
$ws.Range("D240").Value = 44754
$ws.Range("E240").Value = 16
$ws.Range("F240").Value = 100114013
$ws.Range("G240").Value = "Zanahoria"
$ws.Range("H240").Value = "Sin especificar"
$ws.Range("I240").Value = "Primera"
$ws.Range("J240").Value = 100
$ws.Range("K240").Value = 7000
$ws.Range("L240").Value = 7500
$ws.Range("M240").Value = 7250
$ws.Range("N240").Value = "$/saco 20 kilos"
$ws.Range("O240").Value = "Provincia de Diguillín"
$ws.Range("P240").Value = 362
$ws.Range("Q240").Value = 20
$ws.Range("R240").Value = "Hortaliza"
